$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet "results": add the new "S*-unmerged" / "S*-unmergedND" columns
# -------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("results")

# Insert a new column before the current "var" column (G) for "S*-unmerged"
$ws1.Range("G1").EntireColumn.Insert()
$ws1.Range("G1").Value = "S*-unmerged"
$ws1.Range("G2").Value = 128
$ws1.Range("H2").Value = 0
$ws1.Range("J2").Value = $false

# Append a new column after the last one (L) for "S*-unmergedND"
$ws1.Range("L1").Copy()
$ws1.Range("M1").PasteSpecial(-4122)
$ws1.Range("M1").Value = "S*-unmergedND"
$ws1.Range("M2").Value = $false

# -------------------------------------------------------------------------
# Sheet "stats": insert a new "S*-unmerged" row into both data blocks
# -------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("stats")

# Break the existing vertical merges before reshaping rows
$ws2.Range("A2:A6").UnMerge()
$ws2.Range("A7:A11").UnMerge()

# Insert rows bottom-up so row numbers above the insertion point stay put
$ws2.Range("A11").EntireRow.Insert()
$ws2.Range("A6").EntireRow.Insert()

# Restore the label-column style on the freshly inserted rows
$ws2.Range("A5").Copy()
$ws2.Range("A6").PasteSpecial(-4122)
$ws2.Range("B5").Copy()
$ws2.Range("B6").PasteSpecial(-4122)
$ws2.Range("A11").Copy()
$ws2.Range("A12").PasteSpecial(-4122)
$ws2.Range("B11").Copy()
$ws2.Range("B12").PasteSpecial(-4122)

# Write the refreshed figures for every row of both blocks
$ws2.Range("A2").Value = "run 0"
$ws2.Range("B2").Value = "S*-BS"
$ws2.Range("C2").Value = 48
$ws2.Range("D2").Value = 0.0000891457311809063
$ws2.Range("E2").Value = 0.02504725521430373
$ws2.Range("F2").Value = 48
$ws2.Range("G2").Value = 0.002023348119109869
$ws2.Range("H2").Value = 0.003185389097779989
$ws2.Range("I2").Value = 0.006114192772656679
$ws2.Range("J2").Value = 0.01127113122493029
$ws2.Range("K2").Value = 0.0007171705365180969

$ws2.Range("B3").Value = "S*-HS"
$ws2.Range("C3").Value = 48
$ws2.Range("D3").Value = 0.002587008755654097
$ws2.Range("E3").Value = 0.03498434787616134
$ws2.Range("F3").Value = 48
$ws2.Range("G3").Value = 0.00231133634224534
$ws2.Range("H3").Value = 0.006467015016824007
$ws2.Range("I3").Value = 0.007225723005831242
$ws2.Range("J3").Value = 0.01598513964563608
$ws2.Range("K3").Value = 0.000906776636838913

$ws2.Range("B4").Value = "S*-MM"
$ws2.Range("C4").Value = 48
$ws2.Range("D4").Value = 0.002929155249148607
$ws2.Range("E4").Value = 0.03561321692541242
$ws2.Range("F4").Value = 48
$ws2.Range("G4").Value = 0.002260776236653328
$ws2.Range("H4").Value = 0.006592562887817621
$ws2.Range("I4").Value = 0.007945331279188395
$ws2.Range("J4").Value = 0.01615219935774803
$ws2.Range("K4").Value = 0.0007207700982689857

$ws2.Range("B5").Value = "S*-MM0"
$ws2.Range("C5").Value = 48
$ws2.Range("D5").Value = 0.0001705302856862545
$ws2.Range("E5").Value = 0.02893384778872132
$ws2.Range("F5").Value = 48
$ws2.Range("G5").Value = 0.002520418725907803
$ws2.Range("H5").Value = 0.004113330971449614
$ws2.Range("I5").Value = 0.007866930682212114
$ws2.Range("J5").Value = 0.01168011547997594
$ws2.Range("K5").Value = 0.000784547533839941

$ws2.Range("B6").Value = "S*-unmerged"
$ws2.Range("C6").Value = 96
$ws2.Range("D6").Value = 0.003886665217578411
$ws2.Range("E6").Value = 0.08259139815345407
$ws2.Range("F6").Value = 96
$ws2.Range("G6").Value = 0.004012365825474262
$ws2.Range("H6").Value = 0.01105505228042603
$ws2.Range("I6").Value = 0.03809424489736557
$ws2.Range("J6").Value = 0.02431661868467927
$ws2.Range("K6").Value = 0.001474921125918627

$ws2.Range("B7").Value = "Kruskal"
$ws2.Range("C7").Value = 1250
$ws2.Range("E7").Value = 0.01371327193919569

$ws2.Range("A8").Value = "Average"
$ws2.Range("B8").Value = "S*-BS"
$ws2.Range("C8").Value = 48
$ws2.Range("D8").Value = 0.0000891457311809063
$ws2.Range("E8").Value = 0.02504725521430373
$ws2.Range("F8").Value = 48
$ws2.Range("G8").Value = 0.002023348119109869
$ws2.Range("H8").Value = 0.003185389097779989
$ws2.Range("I8").Value = 0.006114192772656679
$ws2.Range("J8").Value = 0.01127113122493029
$ws2.Range("K8").Value = 0.0007171705365180969

$ws2.Range("B9").Value = "S*-HS"
$ws2.Range("C9").Value = 48
$ws2.Range("D9").Value = 0.002587008755654097
$ws2.Range("E9").Value = 0.03498434787616134
$ws2.Range("F9").Value = 48
$ws2.Range("G9").Value = 0.00231133634224534
$ws2.Range("H9").Value = 0.006467015016824007
$ws2.Range("I9").Value = 0.007225723005831242
$ws2.Range("J9").Value = 0.01598513964563608
$ws2.Range("K9").Value = 0.000906776636838913

$ws2.Range("B10").Value = "S*-MM"
$ws2.Range("C10").Value = 48
$ws2.Range("D10").Value = 0.002929155249148607
$ws2.Range("E10").Value = 0.03561321692541242
$ws2.Range("F10").Value = 48
$ws2.Range("G10").Value = 0.002260776236653328
$ws2.Range("H10").Value = 0.006592562887817621
$ws2.Range("I10").Value = 0.007945331279188395
$ws2.Range("J10").Value = 0.01615219935774803
$ws2.Range("K10").Value = 0.0007207700982689857

$ws2.Range("B11").Value = "S*-MM0"
$ws2.Range("C11").Value = 48
$ws2.Range("D11").Value = 0.0001705302856862545
$ws2.Range("E11").Value = 0.02893384778872132
$ws2.Range("F11").Value = 48
$ws2.Range("G11").Value = 0.002520418725907803
$ws2.Range("H11").Value = 0.004113330971449614
$ws2.Range("I11").Value = 0.007866930682212114
$ws2.Range("J11").Value = 0.01168011547997594
$ws2.Range("K11").Value = 0.000784547533839941

$ws2.Range("B12").Value = "S*-unmerged"
$ws2.Range("C12").Value = 96
$ws2.Range("D12").Value = 0.003886665217578411
$ws2.Range("E12").Value = 0.08259139815345407
$ws2.Range("F12").Value = 96
$ws2.Range("G12").Value = 0.004012365825474262
$ws2.Range("H12").Value = 0.01105505228042603
$ws2.Range("I12").Value = 0.03809424489736557
$ws2.Range("J12").Value = 0.02431661868467927
$ws2.Range("K12").Value = 0.001474921125918627

$ws2.Range("B13").Value = "Kruskal"
$ws2.Range("C13").Value = 1250
$ws2.Range("E13").Value = 0.01371327193919569

# Re-establish the vertical merges over the now-6-row blocks
$ws2.Range("A2:A7").Merge()
$ws2.Range("A8:A13").Merge()
